# ============================================================================
# finish inital calculations and documentation
# - adds 4 new worksheets: biodiversity, nitrogen_phosphorus_cycles,
#   freshwater, chemical_pollution
# - wires new D6/D7/row11/row12 cells into the "climate" sheet
# - updates sheet selections / active sheet
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1. climate sheet (sheet2) tweaks
# ----------------------------------------------------------------------------
$climate = $wb.Worksheets.Item("climate")

$climate.Cells.Item(6, 4).Value = "15$/tCO2"          # D6
$climate.Cells.Item(7, 4).Formula = "=B2*F8"           # D7 (was =B2*15*0.74)

$climate.Cells.Item(11, 2).Value = "low"               # B11
$climate.Cells.Item(11, 3).Value = "mid"               # C11
$climate.Cells.Item(11, 4).Value = "high"              # D11

$climate.Cells.Item(12, 2).Formula = "=B8"             # B12
$climate.Cells.Item(12, 3).Formula = "=D7"             # C12
$climate.Cells.Item(12, 4).Formula = "=C9"             # D12

# ----------------------------------------------------------------------------
# 2. new sheet: biodiversity
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$biodiversity = $wb.Worksheets.Add($null, $lastSheet)
$biodiversity.Name = "biodiversity"

$biodiversity.Cells.Item(1, 2).Value = "low"
$biodiversity.Cells.Item(1, 3).Value = "mid"
$biodiversity.Cells.Item(1, 4).Value = "high"

$biodiversity.Cells.Item(2, 1).Value = "global shortfall (USD)"
$biodiversity.Cells.Item(2, 2).Value = 3400000000
$biodiversity.Cells.Item(2, 3).Value = 12700000000
$biodiversity.Cells.Item(2, 4).Value = 38100000000
$biodiversity.Cells.Item(2, 4).NumberFormat = "0.00E+00"

$biodiversity.Cells.Item(3, 1).Value = "per uk person (GBP)"
$biodiversity.Cells.Item(3, 2).Formula = "=B2*Forrest_loss!`$B5*0.74"
$biodiversity.Cells.Item(3, 3).Formula = "=C2*Forrest_loss!`$B5*0.74"
$biodiversity.Cells.Item(3, 4).Formula = "=D2*Forrest_loss!`$B5*0.74"

$biodiversity.Columns.Item(2).ColumnWidth = 11

# ----------------------------------------------------------------------------
# 3. new sheet: nitrogen_phosphorus_cycles
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$nitrogen = $wb.Worksheets.Add($null, $lastSheet)
$nitrogen.Name = "nitrogen_phosphorus_cycles"

$nitrogen.Cells.Item(1, 1).Value = "Fertilizer market size"
$nitrogen.Cells.Item(1, 2).Value = "average profit margins"
$nitrogen.Cells.Item(1, 3).Value = "runoff-rate"
$nitrogen.Cells.Item(1, 4).Value = "cost to reduce"
$nitrogen.Cells.Item(1, 5).Value = "per capita in gbp"

$nitrogen.Cells.Item(2, 1).Value = 158000000000
$nitrogen.Cells.Item(2, 2).Value = 0.3
$nitrogen.Cells.Item(2, 3).Value = 0.5
$nitrogen.Cells.Item(2, 4).Formula = "=C2*A2*B2"
$nitrogen.Cells.Item(2, 5).Formula = "=D2*Forrest_loss!B5*0.74"

$nitrogen.Cells.Item(4, 1).Value = "Arable land for livestock feed"
$nitrogen.Cells.Item(4, 2).Value = "runoff fertilizer reduction from switching to vegan"
$nitrogen.Cells.Item(4, 3).Value = "per capita in gbp"

$nitrogen.Cells.Item(5, 1).Value = 0.36
$nitrogen.Cells.Item(5, 2).Value = 0.7
$nitrogen.Cells.Item(5, 3).Value = 2
$nitrogen.Cells.Item(5, 4).Formula = "=1/(A5*B5*C5/7600000000)"
$nitrogen.Cells.Item(5, 5).Formula = "=D5*Forrest_loss!B5*0.74"

$nitrogen.Cells.Item(7, 1).Value = "Meat consumption rate UK vs global"
$nitrogen.Cells.Item(7, 2).Value = "Low"
$nitrogen.Cells.Item(7, 3).Value = "High"

$nitrogen.Cells.Item(8, 1).Formula = "=B8/2"
$nitrogen.Cells.Item(8, 2).Formula = "=E5"
$nitrogen.Cells.Item(8, 3).Formula = "=E2*2"

$nitrogen.Columns.Item(4).ColumnWidth = 9.140625

# ----------------------------------------------------------------------------
# 4. new sheet: freshwater
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$freshwater = $wb.Worksheets.Add($null, $lastSheet)
$freshwater.Name = "freshwater"

$freshwater.Cells.Item(1, 1).Value = "agriculture water use rate"
$freshwater.Cells.Item(1, 2).Value = "of which for livestock feed"
$freshwater.Cells.Item(1, 3).Value = "reduction by going vegan"
$freshwater.Cells.Item(1, 4).Value = "per capita in gbp"

$freshwater.Cells.Item(2, 1).Value = 0.92
$freshwater.Cells.Item(2, 2).Value = 0.29
$freshwater.Cells.Item(2, 3).Value = 0.8
$freshwater.Cells.Item(2, 4).Value = 2
$freshwater.Cells.Item(2, 5).Formula = "=1/(A2*C2*B2*D2/7600000000)"
$freshwater.Cells.Item(2, 6).Formula = "=E2*Forrest_loss!B5*0.74"

$freshwater.Cells.Item(6, 2).Value = "low"
$freshwater.Cells.Item(6, 3).Value = "mid"
$freshwater.Cells.Item(6, 4).Value = "high"

$freshwater.Cells.Item(7, 2).Formula = "=C7/2"
$freshwater.Cells.Item(7, 3).Formula = "=D7/2"
$freshwater.Cells.Item(7, 4).Formula = "=F2"

$freshwater.Columns.Item(5).ColumnWidth = 10.140625

# ----------------------------------------------------------------------------
# 5. new sheet: chemical_pollution
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$chem = $wb.Worksheets.Add($null, $lastSheet)
$chem.Name = "chemical_pollution"

$chem.Cells.Item(1, 1).Value = "Plastics"

$chem.Cells.Item(2, 3).Value = "$/kg"
$chem.Cells.Item(2, 4).Value = "overall cost"

$chem.Cells.Item(3, 2).Value = "Annual waste production"
$chem.Cells.Item(3, 3).Value = 415000000000
$chem.Cells.Item(3, 4).Formula = "=0.21*365"
$chem.Cells.Item(3, 5).Value = "tCO2/kg"
$chem.Cells.Item(3, 8).Value = "15$/tCO2"
$chem.Cells.Item(3, 9).Value = "cost per kg"

$chem.Cells.Item(4, 2).Value = "Going to landfills"
$chem.Cells.Item(4, 3).Value = 201000000000
$chem.Cells.Item(4, 5).Value = "equal proportion"
$chem.Cells.Item(4, 6).Value = 0.0025
$chem.Cells.Item(4, 7).Value = "15$/tCO2"
$chem.Cells.Item(4, 8).Value = 5
$chem.Cells.Item(4, 9).Formula = "=I5/C4"
$chem.Cells.Item(4, 10).Formula = "=J5/C4"

$chem.Cells.Item(5, 2).Value = "Going to mismanaged coastal"
$chem.Cells.Item(5, 3).Value = 37600000000
$chem.Cells.Item(5, 4).Formula = "=C5*(0.0021/3+0.0014*2/3)/68000000"
$chem.Cells.Item(5, 5).Value = "capital cost per tonne of annual capacity"
$chem.Cells.Item(5, 6).Value = 1000
$chem.Cells.Item(5, 7).Value = "overall cost"
$chem.Cells.Item(5, 8).Formula = "=C6*H4"
$chem.Cells.Item(5, 9).Value = 660000000000
$chem.Cells.Item(5, 10).Value = 94000000000

$chem.Cells.Item(6, 2).Value = "Going into the ocean"
$chem.Cells.Item(6, 3).Value = 9400000000
$chem.Cells.Item(6, 5).Value = "cost per kg"
$chem.Cells.Item(6, 6).Formula = "=F4*climate!F8+F5/20/1000"

$chem.Cells.Item(8, 2).Value = "Incineration"
$chem.Cells.Item(8, 6).Formula = "=0.74*F6/7600000000*C4"
$chem.Cells.Item(8, 8).Formula = "=0.74*H5/7600000000"
$chem.Cells.Item(8, 9).Formula = "=0.74*I5/7600000000"
$chem.Cells.Item(8, 10).Formula = "=0.74*J5/7600000000"

$chem.Cells.Item(9, 2).Value = "OceanCleanup"
$chem.Cells.Item(9, 6).Formula = "=0.74*F6*Forrest_loss!B5*C4"
$chem.Cells.Item(9, 8).Formula = "=0.74*H5*Forrest_loss!`$B5"
$chem.Cells.Item(9, 9).Formula = "=0.74*I5*Forrest_loss!`$B5"
$chem.Cells.Item(9, 10).Formula = "=0.74*J5*Forrest_loss!`$B5"

$chem.Cells.Item(10, 2).Value = "world"
$chem.Cells.Item(10, 3).Value = 0.02
$chem.Cells.Item(10, 6).Formula = "=((`$D`$5*7600000000/`$C`$5+`$C10)*`$D`$3)*F`$6"
$chem.Cells.Item(10, 8).Formula = "=(`$C10*(`$C`$5/`$C`$3*`$D`$3)+`$D`$5)*H`$4*0.74"
$chem.Cells.Item(10, 9).Formula = "=(`$C10*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*I`$4*0.74"
$chem.Cells.Item(10, 10).Formula = "=(`$C10*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*J`$4*0.74"

$chem.Cells.Item(11, 3).Value = 0.1
$chem.Cells.Item(11, 6).Formula = "=((`$D`$5*7600000000/`$C`$5+`$C11)*`$D`$3)*F`$6"
$chem.Cells.Item(11, 8).Formula = "=(`$C11*(`$C`$5/`$C`$3*`$D`$3)+`$D`$5)*H`$4*0.74"
$chem.Cells.Item(11, 9).Formula = "=(`$C11*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*I`$4*0.74"
$chem.Cells.Item(11, 10).Formula = "=(`$C11*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*J`$4*0.74"

$chem.Cells.Item(12, 3).Value = 0.2
$chem.Cells.Item(12, 6).Formula = "=((`$D`$5*7600000000/`$C`$5+`$C12)*`$D`$3)*F`$6"
$chem.Cells.Item(12, 8).Formula = "=(`$C12*(`$C`$5/`$C`$3*`$D`$3)+`$D`$5)*H`$4*0.74"
$chem.Cells.Item(12, 9).Formula = "=(`$C12*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*I`$4*0.74"
$chem.Cells.Item(12, 10).Formula = "=(`$C12*(`$C`$4/`$C`$3*`$D`$3)+`$D`$5*`$C`$4/`$C`$5)*J`$4*0.74"

$chem.Cells.Item(14, 6).Value = "low"
$chem.Cells.Item(14, 7).Value = "med"
$chem.Cells.Item(14, 8).Value = "high"

$chem.Cells.Item(15, 6).Formula = "=AVERAGE(H10,J10)"
$chem.Cells.Item(15, 7).Formula = "=AVERAGE(H12,J12)"
$chem.Cells.Item(15, 8).Formula = "=AVERAGE(H9,J9)"

$chem.Cells.Item(17, 1).Value = "POPs"

$chem.Cells.Item(18, 2).Value = "low est"

$chem.Cells.Item(19, 2).Value = "DDT phase out cost to developing nations"
$chem.Cells.Item(19, 3).Value = "high est"

$chem.Cells.Item(20, 2).Value = 350000000
$chem.Cells.Item(20, 2).NumberFormat = "0.00E+00"
$chem.Cells.Item(20, 3).Value = 950000000
$chem.Cells.Item(20, 3).NumberFormat = "0.00E+00"

$chem.Cells.Item(22, 6).Value = "low"
$chem.Cells.Item(22, 7).Value = "med"
$chem.Cells.Item(22, 8).Value = "high"

$chem.Cells.Item(23, 6).Formula = "=B20*Forrest_loss!`$B5*0.74"
$chem.Cells.Item(23, 6).NumberFormat = "0.00E+00"
$chem.Cells.Item(23, 7).Formula = "=C20*Forrest_loss!`$B5*0.74"
$chem.Cells.Item(23, 7).NumberFormat = "0.00E+00"
$chem.Cells.Item(23, 8).Formula = "=C20*Forrest_loss!`$B5*0.74*2"
$chem.Cells.Item(23, 8).NumberFormat = "0.00E+00"

$chem.Cells.Item(25, 1).Value = "Others"

$chem.Cells.Item(26, 6).Value = "low"
$chem.Cells.Item(26, 7).Value = "med"
$chem.Cells.Item(26, 8).Value = "high"

$chem.Cells.Item(27, 6).Formula = "=1000000000*Forrest_loss!B5*0.74"
$chem.Cells.Item(27, 7).Formula = "=10000000000*Forrest_loss!B5*0.74"
$chem.Cells.Item(27, 8).Formula = "=50000000000*Forrest_loss!B5*0.74"

$chem.Columns.Item(1).ColumnWidth = 9.7109375
$chem.Columns.Item(2).ColumnWidth = 18.85546875

# ----------------------------------------------------------------------------
# 6. selections per-sheet (also fixes which sheet / cell is active on open)
# ----------------------------------------------------------------------------
$forrest = $wb.Worksheets.Item("Forrest_loss")
$forrest.Activate()
$forrest.Range("E6").Select()

$climate.Activate()
$climate.Range("D7").Select()

$biodiversity.Activate()
$biodiversity.Range("B3:D3").Select()

$nitrogen.Activate()
$nitrogen.Range("E5").Select()

$freshwater.Activate()
$freshwater.Range("B7:D7").Select()

$chem.Activate()
$chem.Range("F27:H27").Select()
